$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.535.32'
$ws.Range("E2").Value = '  +2.05%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.363.85'
$ws.Range("E3").Value = '  +6.13%  '

# Row 4
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.96'
$ws.Range("E5").Value = '  +6.65%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.66'
$ws.Range("E6").Value = '  -3.15%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.640'
$ws.Range("E7").Value = '  +3.08%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.638'
$ws.Range("E9").Value = '  +4.88%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.68'
$ws.Range("E10").Value = '  -3.78%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0938'
$ws.Range("E11").Value = '  +2.25%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.74'
$ws.Range("E12").Value = '  -1.65%  '

# Row 13
$ws.Range("E13").Value = '  +3.70%  '

# Row 14
$ws.Range("E14").Value = '  +2.50%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.62'
$ws.Range("E15").Value = '  +9.90%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.721.41'
$ws.Range("E16").Value = '  +6.31%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.432.19'
$ws.Range("E17").Value = '  +8.90%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.537.38'
$ws.Range("E18").Value = '  +2.50%  '

# Row 19
$ws.Range("E19").Value = '  +2.67%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.28'
$ws.Range("E20").Value = '  -1.34%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.26'
$ws.Range("E21").Value = '  +3.25%  '

# Row 22
$ws.Range("E22").Value = '  -1.21%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.53'
$ws.Range("E23").Value = '  +7.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '258.74'
$ws.Range("E24").Value = '  +12.68%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.41'
$ws.Range("E25").Value = '  +1.80%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.05'
$ws.Range("E26").Value = '  +2.62%  '

# Row 27
$ws.Range("E27").Value = '  +0.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '39.04'
$ws.Range("E28").Value = '  +1.34%  '

# Row 29
$ws.Range("E29").Value = '  +0.18%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.78'
$ws.Range("E30").Value = '  +7.90%  '

# Row 31
$ws.Range("E31").Value = '  -1.01%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.32'
$ws.Range("E32").Value = '  -0.28%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0924'
$ws.Range("E33").Value = '  +2.09%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.98'
$ws.Range("E34").Value = '  +5.17%  '

# Row 35
$ws.Range("E35").Value = '  +4.84%  '

# Row 36
$ws.Range("E36").Value = '  -4.71%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0374'
$ws.Range("E37").Value = '  -0.94%  '

# Row 38
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.10'
$ws.Range("E38").Value = '  -5.00%  '

# Row 39
$ws.Range("E39").Value = '  +0.16%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.79'
$ws.Range("E40").Value = '  +14.90%  '

# Row 41
$ws.Range("E41").Value = '  +13.61%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.95'
$ws.Range("E42").Value = '  -0.30%  '

# Row 43
$ws.Range("E43").Value = '  -1.19%  '

# Row 44
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.71'
$ws.Range("E45").Value = '  -0.79%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.63'
$ws.Range("E46").Value = '  +2.51%  '

# Row 47
$ws.Range("E47").Value = '  +8.81%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '111.72'
$ws.Range("E48").Value = '  +7.98%  '

# Row 49
$ws.Range("E49").Value = '  -1.65%  '

# Row 50
$ws.Range("E50").Value = '  +2.58%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.474'
$ws.Range("E51").Value = '  +7.36%  '
